# FlixelRL-199 Floor4 drop table
# Rework the "item_appear" sheet: rows 3-22 get new id/ratio values (Floor4
# drop weights), and four new rows (23-26) are appended for the remaining
# Portion tiers. Column A/C get new widths too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("item_appear")

# --- column widths -----------------------------------------------------
# (ColumnWidth is quantised to the host's fixed digit-width grid, so these
# are the closest attainable values to the authored widths of 2.5 / 5.90625)
$ws.Columns.Item(1).ColumnWidth = 1.7142857142857142
$ws.Columns.Item(3).ColumnWidth = 5.142857142857143

# --- existing rows 3-22 -------------------------------------------------
$rows = @(
    @{ Row=3;  A=1;  B=1; C=10; D="WEAPON1";  E=50 }
    @{ Row=4;  A=2;  B=1; C=10; D="WEAPON2";  E=30 }
    @{ Row=5;  A=3;  B=1; C=10; D="WEAPON3";  E=20 }
    @{ Row=6;  A=4;  B=1; C=10; D="WEAPON4";  E=10 }
    @{ Row=7;  A=5;  B=1; C=10; D="WEAPON5";  E=5 }
    @{ Row=8;  A=6;  B=1; C=10; D="WEAPON6";  E=2 }
    @{ Row=9;  A=7;  B=1; C=10; D="WEAPON7";  E=1 }
    @{ Row=10; A=8;  B=1; C=10; D="ARMOR1";   E=50 }
    @{ Row=11; A=9;  B=1; C=10; D="ARMOR2";   E=30 }
    @{ Row=12; A=10; B=1; C=10; D="ARMOR3";   E=20 }
    @{ Row=13; A=11; B=1; C=10; D="ARMOR4";   E=10 }
    @{ Row=14; A=12; B=1; C=10; D="ARMOR5";   E=5 }
    @{ Row=15; A=13; B=1; C=10; D="ARMOR6";   E=2 }
    @{ Row=16; A=14; B=1; C=10; D="ARMOR7";   E=1 }
    @{ Row=17; A=15; B=1; C=10; D="FOOD1";    E=200 }
    @{ Row=18; A=16; B=1; C=10; D="FOOD2";    E=50 }
    @{ Row=19; A=17; B=1; C=10; D="FOOD3";    E=20 }
    @{ Row=20; A=18; B=1; C=10; D="FOOD4";    E=20 }
    @{ Row=21; A=19; B=1; C=10; D="PORTION1"; E=500 }
    @{ Row=22; A=20; B=1; C=10; D="PORTION2"; E=400 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value = $r.A
    $ws.Cells.Item($n, 2).Value = $r.B
    $ws.Cells.Item($n, 3).Value = $r.C
    $ws.Cells.Item($n, 4).Value = $r.D
    $ws.Cells.Item($n, 5).Value = $r.E
}

# --- new rows 23-26 -------------------------------------------------
# Clone row 22's formatting (style + row height) down into the new rows
# before stamping in their values, so they match the rest of the table.
$ws.Range("A22:E22").Copy($ws.Range("A23:E26"))
$ws.Rows.Item(23).RowHeight = 20
$ws.Rows.Item(24).RowHeight = 20
$ws.Rows.Item(25).RowHeight = 20
$ws.Rows.Item(26).RowHeight = 20

$newRows = @(
    @{ Row=23; A=21; B=1; C=10; D="PORTION3"; E=300 }
    @{ Row=24; A=22; B=1; C=10; D="PORTION4"; E=100 }
    @{ Row=25; A=23; B=1; C=10; D="PORTION5"; E=50 }
    @{ Row=26; A=24; B=1; C=10; D="PORTION6"; E=25 }
)

foreach ($r in $newRows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value = $r.A
    $ws.Cells.Item($n, 2).Value = $r.B
    $ws.Cells.Item($n, 3).Value = $r.C
    $ws.Cells.Item($n, 4).Value = $r.D
    $ws.Cells.Item($n, 5).Value = $r.E
}
